$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = " Ajmer"
$ws.Range("B2").Value = "Ajmer"
$ws.Range("A3").Value = " Alwar"
$ws.Range("B3").Value = "Alwar"
$ws.Range("A4").Value = " Bhilai"
$ws.Range("B4").Value = "Bhilai"

$ws.Range("B4").Select()
